# Update cryptos list (Price / Volume(1h) columns) with refreshed quote data.
# Values that look like plain numbers are prefixed with a leading apostrophe
# so Excel stores them as text (matching the source data's inline-string
# format, e.g. "42.80") instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.078.90"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "2.219.51"
$ws.Range("E3").Value = "  -1.21%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'241.47"
$ws.Range("E5").Value = "  -2.29%  "
$ws.Range("D6").Value = "'0.627"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("D7").Value = "'73.52"
$ws.Range("E7").Value = "  -1.29%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "'0.608"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").Value = "'42.80"
$ws.Range("E10").Value = "  +1.42%  "
$ws.Range("D11").Value = "'0.0957"
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").Value = "'7.09"
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("D13").Value = "'0.104"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").Value = "2.552.31"
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("D15").Value = "'14.33"
$ws.Range("E15").Value = "  -1.78%  "
$ws.Range("D16").Value = "'0.839"
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("D17").Value = "2.213.96"
$ws.Range("E17").Value = "  -2.65%  "
$ws.Range("D18").Value = "41.925.72"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").Value = "'0.0000107"
$ws.Range("E19").Value = "  +8.37%  "
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("D21").Value = "'6.18"
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("D22").Value = "'10.76"
$ws.Range("E22").Value = "  +19.88%  "
$ws.Range("D23").Value = "'229.89"
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("D24").Value = "'2.11"
$ws.Range("E24").Value = "  -5.54%  "
$ws.Range("D25").Value = "'11.80"
$ws.Range("E25").Value = "  +3.01%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").Value = "'2.27"
$ws.Range("E28").Value = "  -2.04%  "
$ws.Range("D29").Value = "'2.20"
$ws.Range("E29").Value = "  +0.70%  "
$ws.Range("D30").Value = "'167.49"
$ws.Range("E30").Value = "  -1.00%  "
$ws.Range("D31").Value = "'20.52"
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("D32").Value = "'5.58"
$ws.Range("E32").Value = "  +7.19%  "
$ws.Range("D33").Value = "'0.0797"
$ws.Range("E33").Value = "  -3.39%  "
$ws.Range("D34").Value = "'29.65"
$ws.Range("E34").Value = "  -3.85%  "
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("E36").Value = "  -9.35%  "
$ws.Range("D37").Value = "'4.27"
$ws.Range("E37").Value = "  -4.42%  "
$ws.Range("D38").Value = "'0.0301"
$ws.Range("E38").Value = "  -5.04%  "
$ws.Range("D39").Value = "'13.65"
$ws.Range("E39").Value = "  -2.03%  "
$ws.Range("D40").Value = "'66.22"
$ws.Range("E40").Value = "  +5.12%  "
$ws.Range("D41").Value = "'2.13"
$ws.Range("E41").Value = "  -2.52%  "
$ws.Range("D42").Value = "'5.64"
$ws.Range("E42").Value = "  -2.56%  "
$ws.Range("D43").Value = "'0.199"
$ws.Range("E43").Value = "  -2.59%  "
$ws.Range("D44").Value = "'8.76"
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("D45").Value = "'105.23"
$ws.Range("E45").Value = "  -1.58%  "
$ws.Range("E46").Value = "  -2.75%  "
$ws.Range("D47").Value = "'2.42"
$ws.Range("E47").Value = "  +4.44%  "
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("D49").Value = "'1.17"
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("D50").Value = "'2.69"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").Value = "2.424.00"
$ws.Range("E51").Value = "  -1.35%  "
